# Update gh-pages output data (stats refresh) across the three data sheets:
# 展览 (Exhibitions), 演出 (Shows), 全部类型 (All types).
# Column F = 想去人数 (interest count), Column G = 最低票价 (min ticket price).

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 451
$ws.Range("G2").Value = 88
$ws.Range("F3").Value = 1838
$ws.Range("F4").Value = 1451
$ws.Range("F6").Value = 1735
$ws.Range("F9").Value = 653
$ws.Range("F10").Value = 30
$ws.Range("F16").Value = 20
$ws.Range("F19").Value = 102
$ws.Range("F20").Value = 4538
$ws.Range("F21").Value = 40
$ws.Range("F24").Value = 2162
$ws.Range("F26").Value = 5
$ws.Range("F27").Value = 2028

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 71

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 451
$ws.Range("G2").Value = 88
$ws.Range("F3").Value = 1838
$ws.Range("F4").Value = 1451
$ws.Range("F6").Value = 1735
$ws.Range("F9").Value = 653
$ws.Range("F10").Value = 30
$ws.Range("F16").Value = 20
$ws.Range("F19").Value = 102
$ws.Range("F20").Value = 4539
$ws.Range("F21").Value = 71
$ws.Range("F22").Value = 40
$ws.Range("F26").Value = 2162
$ws.Range("F28").Value = 5
$ws.Range("F29").Value = 2028

Write-Output "Applied gh-pages data refresh to 展览, 演出, 全部类型 sheets"
